$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.0883
$ws.Range("A9").Value = -21.9657
$ws.Range("E12").Value = 18.05830000000002
$ws.Range("A18").Value = -22.34610000000001
$ws.Range("A20").Value = -20.55309999999998
$ws.Range("E26").Value = 16.1223
$ws.Range("A27").Value = -22.0963
$ws.Range("E27").Value = 16.68479999999998
$ws.Range("E29").Value = 16.94150000000001
$ws.Range("E37").Value = 16.74050000000001
$ws.Range("E38").Value = 16.431
$ws.Range("E51").Value = 17.30840000000001
$ws.Range("E55").Value = 16.51030000000001
$ws.Range("A69").Value = -21.68309999999998
$ws.Range("E69").Value = 17.38780000000002
$ws.Range("E70").Value = 18.08940000000002
$ws.Range("A76").Value = -19.99319999999998
$ws.Range("A82").Value = -22.0505
$ws.Range("E83").Value = 16.55099999999999
$ws.Range("E102").Value = 16.8121
